# Adds the "adjusted_2020" column (inflation-adjusted spending) to Sheet1,
# formats it like currency with a dark-red font, resizes the column and
# updates the sheet selection, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column E.
$ws.Range("E1").Value = "adjusted_2020"

# Inflation-adjusted ("adjusted_2020") values for each year row (2-26).
$adjusted = @(
    12092.88,
    11375.12,
    13118.42,
    13056.59,
    15324.55,
    16357.01,
    15324.16,
    15704.62,
    24211.79,
    25847.25,
    26895.13,
    26816.89,
    27634.87,
    26684.62,
    30379.15,
    29311.64,
    28541.83,
    27888.5,
    26672.69,
    28379.39,
    28581.16,
    29262.32,
    30608.74,
    34643.38,
    37723.56
)

for ($i = 0; $i -lt $adjusted.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $adjusted[$i]
}

# Style the new data like currency: dark red font + 2-decimal thousands
# format. Build the style once on E2, then fan it out to the rest of the
# column via copy/paste-special (format only) so every cell shares the same
# single cellXf instead of each property-assignment minting its own.
$ws.Range("E2").Font.Color = 51
$ws.Range("E2").NumberFormat = "#,##0.00"

$ws.Range("E2").Copy()
$ws.Range("E3:E26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fit the new column's width.
$ws.Columns.Item(5).ColumnWidth = 11

# Update the view: clear the frozen/scrolled top-left cell and select E11.
[void]$ws.Range("E11").Select()
